# Slide 1 ("Modules and packages"): the "Content Placeholder 2" shape
# (shape id=3) gets a new final paragraph reading "heapq" appended after
# the existing "__name__, '__main__'" paragraph.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$target = $s.Shapes.Item("Content Placeholder 2")

$tr = $target.TextFrame.TextRange
# vbCr (carriage return) starts a new paragraph; InsertAfter appends after
# the current end of the text range, preserving all existing runs/formatting.
[void]$tr.InsertAfter([char]13 + "heapq")
